# Natmi following Dr Hou advice
#
# The sending-cluster set grows from {FAPs, sCs} to {ECs, FAPs, sCs}, so the
# sheet now has the full 3x3 (sending cluster x target cluster) cross join
# for the Il6 -> Il6st ligand/receptor pair (9 data rows instead of 6), and
# every previously-existing row's metrics are recomputed against the new
# 3-cluster population.
#
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E Ligand-expressing cells, F Ligand detection
#          rate, G Ligand average expression value, H Ligand total
#          expression value, I/J Ligand derived specificity (avg/total),
#          K Receptor-expressing cells, L Receptor detection rate,
#          M Receptor average expression value, N Receptor total expression
#          value, O/P Receptor derived specificity (avg/total), Q/R Edge
#          average/total expression weight, S/T Edge average/total
#          expression derived specificity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Il6", "Il6st", "ECs",  2, 0.6666666666666666, 5.67446,           17.02338,          0.245923038656222, 0.245923038656222, 3, 1, 37.42645866666667, 112.279376,        0.2415534622699011, 0.2415534622699011, 212.3749426456533, 1911.37448381088,  0.05940356143934516, 0.05940356143934516),
    @("ECs",  "Il6", "Il6st", "FAPs", 2, 0.6666666666666666, 5.67446,           17.02338,          0.245923038656222, 0.245923038656222, 3, 1, 94.96115633333334, 284.883469,        0.6128871635375853, 0.6128871635375853, 538.8532831672467, 4849.67954850522,  0.1507230736105559,  0.1507230736105559),
    @("ECs",  "Il6", "Il6st", "sCs",  2, 0.6666666666666666, 5.67446,           17.02338,          0.245923038656222, 0.245923038656222, 3, 1, 22.553069,         67.659207,         0.1455593741925136, 0.1455593741925136, 127.97648791774,   1151.78839125966,  0.035796403606321,   0.035796403606321),
    @("FAPs", "Il6", "Il6st", "ECs",  3, 1,                  14.93150033333333, 44.794501,         0.647110021688359, 0.647110021688359, 3, 1, 37.42645866666667, 112.279376,        0.2415534622699011, 0.2415534622699011, 558.8331800568195, 5029.498620511376, 0.1563116662083739,  0.1563116662083739),
    @("FAPs", "Il6", "Il6st", "FAPs", 3, 1,                  14.93150033333333, 44.794501,         0.647110021688359, 0.647110021688359, 3, 1, 94.96115633333334, 284.883469,        0.6128871635375853, 0.6128871635375853, 1417.912537444885, 12761.21283700397, 0.3966054256893237,  0.3966054256893237),
    @("FAPs", "Il6", "Il6st", "sCs",  3, 1,                  14.93150033333333, 44.794501,         0.647110021688359, 0.647110021688359, 3, 1, 22.553069,         67.659207,         0.1455593741925136, 0.1455593741925136, 336.7511572911896, 3030.760415620707, 0.0941929297906614,  0.0941929297906614),
    @("sCs",  "Il6", "Il6st", "ECs",  3, 1,                  2.468169,          7.404507,          0.106966939655419, 0.106966939655419, 3, 1, 37.42645866666667, 112.279376,        0.2415534622699011, 0.2415534622699011, 92.37482506084801, 831.373425547632,  0.02583823462218204, 0.02583823462218204),
    @("sCs",  "Il6", "Il6st", "FAPs", 3, 1,                  2.468169,          7.404507,          0.106966939655419, 0.106966939655419, 3, 1, 94.96115633333334, 284.883469,        0.6128871635375853, 0.6128871635375853, 234.380182266087,  2109.421640394783, 0.06555866423770582, 0.0655586642377058),
    @("sCs",  "Il6", "Il6st", "sCs",  3, 1,                  2.468169,          7.404507,          0.106966939655419, 0.106966939655419, 3, 1, 22.553069,         67.659207,         0.1455593741925136, 0.1455593741925136, 55.664785760661,   500.9830718459489, 0.01557004079553115, 0.01557004079553115)
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

Write-Host "applied $($data.Length) rows (A2:T10)"
